# Append the new daily row (2025/10/05) to the tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A holds a literal "yyyy/mm/dd" text label (not a real date), matching
# the rest of the sheet. Excel's input parser would otherwise autoconvert the
# string into a date serial, so force Text format before typing it in, then
# clear the (now unneeded) explicit formatting back to the sheet's default so
# the new row's cells come out unstyled just like every other data row.
$ws.Range("A62").NumberFormat = "@"
$ws.Range("A62").Value = "2025/10/05"
$ws.Range("A62").ClearFormats()

$ws.Range("B62").Value = "日"
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 201
